# Apply updated crypto price/volume data per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.586.51"
$ws.Range("D3").Value = "1.662.72"
$ws.Range("E3").Value = "  -4.07%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'215.41"
$ws.Range("E5").Value = "  -1.93%  "
$ws.Range("D6").Value = "'0.509"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'24.35"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  -2.65%  "
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").Value = "'0.0878"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").Value = "1.899.55"
$ws.Range("E12").Value = "  -4.01%  "
$ws.Range("D13").Value = "1.690.90"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "'66.14"
$ws.Range("D17").Value = "27.578.47"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").Value = "'241.35"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "0.0₃0732"
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("D20").Value = "'7.64"
$ws.Range("E20").Value = "  -3.93%  "
$ws.Range("D22").Value = "'4.51"
$ws.Range("E22").Value = "  -3.31%  "
$ws.Range("E23").Value = "  -3.62%  "
$ws.Range("D24").Value = "'2.06"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").Value = "'146.30"
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D26").Value = "'7.24"
$ws.Range("E26").Value = "  -4.15%  "
$ws.Range("D27").Value = "'16.35"
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("D30").Value = "'1.21"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("E32").Value = "  -2.62%  "
$ws.Range("D33").Value = "1.462.05"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("D34").Value = "'3.12"
$ws.Range("E34").Value = "  -4.51%  "
$ws.Range("D35").Value = "'1.58"
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.38"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'0.929"
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("D39").Value = "'0.572"
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").Value = "'70.18"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("E41").Value = "  -5.18%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  -4.04%  "
$ws.Range("E44").Value = "  -3.64%  "
$ws.Range("D45").Value = "'0.793"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").Value = "1.806.83"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").Value = "'88.79"
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("E49").Value = "  -5.09%  "
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("D51").Value = "'7.90"
$ws.Range("E51").Value = "  -3.63%  "
